$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
function Add-Hyperlink($row, $col, $url, $styleSourceAddr) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = $url
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
    $ws.Range($styleSourceAddr).Copy()
    $cell.PasteSpecial(-4122)
}
Add-Hyperlink 377 23 "https://www.openstreetmap.org/way/1" "W99"
Add-Hyperlink 377 24 "https://www.openstreetmap.org/way/2" "X99"
Add-Hyperlink 377 25 "https://www.openstreetmap.org/way/3" "Y99"
Add-Hyperlink 378 23 "https://www.openstreetmap.org/way/4" "W99"
Write-Host "count:" $ws.Hyperlinks.Count
